$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of column J (budget ratio) with bold font
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true
$ws.Range("J12").Font.Size = 11

# Row 14-17: summary labels and statistics
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$summaryRange = $ws.Range("B14:B17")
$summaryRange.Font.Bold = $true
$summaryRange.Font.Size = 12
$summaryRange.VerticalAlignment = -4108

$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

$ws.Range("A14:B17").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$wb.Save()
